$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2024-09-08 Sunday", $true, $false, $false, $false, $false, $true, 0, $false, "2024-09-09 Monday", 1) | Out-Null

# Update table cell values (addressed by row/col, using wdReplaceOne scoped
# to each cells own Range, so duplicate old text in sibling cells cannot
# leak across cells during replacement)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("25÷3=8, 1", $true, $false, $false, $false, $false, $true, 0, $false, "40÷6=6, 4", 1) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("35÷9=3, 8", $true, $false, $false, $false, $false, $true, 0, $false, "65÷3=21, 2", 1) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("55÷7=7, 6", $true, $false, $false, $false, $false, $true, 0, $false, "59÷5=11, 4", 1) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("34÷4=8, 2", $true, $false, $false, $false, $false, $true, 0, $false, "53÷7=7, 4", 1) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("41÷2=20, 1", $true, $false, $false, $false, $false, $true, 0, $false, "23÷3=7, 2", 1) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("48÷4=12, 0", $true, $false, $false, $false, $false, $true, 0, $false, "18÷6=3, 0", 1) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 0, $false, "92÷3=30, 2", 1) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 0, $false, "23÷9=2, 5", 1) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("28÷9=3, 1", $true, $false, $false, $false, $false, $true, 0, $false, "16÷6=2, 4", 1) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("28÷9=3, 1", $true, $false, $false, $false, $false, $true, 0, $false, "55÷8=6, 7", 1) | Out-Null
$t.Cell(9, 1).Range.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 0, $false, "57÷9=6, 3", 1) | Out-Null
$t.Cell(9, 2).Range.Find.Execute("67÷6=11, 1", $true, $false, $false, $false, $false, $true, 0, $false, "40÷2=20, 0", 1) | Out-Null
$t.Cell(9, 3).Range.Find.Execute("19÷5=3, 4", $true, $false, $false, $false, $false, $true, 0, $false, "99÷3=33, 0", 1) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("66÷3=22, 0", $true, $false, $false, $false, $false, $true, 0, $false, "45÷3=15, 0", 1) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("89÷2=44, 1", $true, $false, $false, $false, $false, $true, 0, $false, "61÷9=6, 7", 1) | Out-Null
$t.Cell(13, 1).Range.Find.Execute("52÷3=17, 1", $true, $false, $false, $false, $false, $true, 0, $false, "48÷4=12, 0", 1) | Out-Null
$t.Cell(13, 2).Range.Find.Execute("52÷4=13, 0", $true, $false, $false, $false, $false, $true, 0, $false, "92÷2=46, 0", 1) | Out-Null
$t.Cell(13, 3).Range.Find.Execute("43÷4=10, 3", $true, $false, $false, $false, $false, $true, 0, $false, "58÷2=29, 0", 1) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("65÷8=8, 1", $true, $false, $false, $false, $false, $true, 0, $false, "35÷5=7, 0", 1) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("77÷6=12, 5", $true, $false, $false, $false, $false, $true, 0, $false, "84÷6=14, 0", 1) | Out-Null
$t.Cell(17, 1).Range.Find.Execute("54÷4=13, 2", $true, $false, $false, $false, $false, $true, 0, $false, "89÷4=22, 1", 1) | Out-Null
$t.Cell(17, 2).Range.Find.Execute("94÷3=31, 1", $true, $false, $false, $false, $false, $true, 0, $false, "17÷5=3, 2", 1) | Out-Null
$t.Cell(17, 3).Range.Find.Execute("14÷5=2, 4", $true, $false, $false, $false, $false, $true, 0, $false, "26÷4=6, 2", 1) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("15÷7=2, 1", $true, $false, $false, $false, $false, $true, 0, $false, "88÷5=17, 3", 1) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("60÷8=7, 4", $true, $false, $false, $false, $false, $true, 0, $false, "25÷9=2, 7", 1) | Out-Null

# Sanity check: verify every targeted cell now holds its expected value.
$expected = @{
    "1,1"="40÷6=6, 4"; "1,2"="65÷3=21, 2"; "1,3"="59÷5=11, 4"; "1,4"="53÷7=7, 4"; "1,5"="23÷3=7, 2";
    "5,1"="18÷6=3, 0"; "5,2"="92÷3=30, 2"; "5,3"="23÷9=2, 5"; "5,4"="16÷6=2, 4"; "5,5"="55÷8=6, 7";
    "9,1"="57÷9=6, 3"; "9,2"="40÷2=20, 0"; "9,3"="99÷3=33, 0"; "9,4"="45÷3=15, 0"; "9,5"="61÷9=6, 7";
    "13,1"="48÷4=12, 0"; "13,2"="92÷2=46, 0"; "13,3"="58÷2=29, 0"; "13,4"="35÷5=7, 0"; "13,5"="84÷6=14, 0";
    "17,1"="89÷4=22, 1"; "17,2"="17÷5=3, 2"; "17,3"="26÷4=6, 2"; "17,4"="88÷5=17, 3"; "17,5"="25÷9=2, 7";
}
foreach ($key in $expected.Keys) {
    $parts = $key.Split(",")
    $row = [int]$parts[0]
    $col = [int]$parts[1]
    $actual = $t.Cell($row, $col).Range.Text.TrimEnd([char]7, [char]13)
    if ($actual -ne $expected[$key]) {
        throw "Cell($row,$col) mismatch: expected '$($expected[$key])' got '$actual'"
    }
}
